# Add: Search a lease
#
# Inserts a new "GlobalSearch" worksheet (object-locator data for the new
# lease global-search UI) right after "Lease" and before "Homepage", fills
# it with the ObjReference/ObjSearchKey/ObjectLocators/ParentLocator table,
# and makes it the active (selected) sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet immediately after "Lease" (pushes Homepage, etc. down
# one slot) and name it.
$lease = $wb.Worksheets.Item("Lease")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lease)
$ws.Name = "GlobalSearch"

# Header row.
$ws.Range("A1").Value = "ObjReference"
$ws.Range("B1").Value = "ObjSearchKey"
$ws.Range("C1").Value = "ObjectLocators"
$ws.Range("D1").Value = "ParentLocator"

# Search link on the menu bar.
$ws.Range("A2").Value = "search"
$ws.Range("B2").Value = "by_xpath"
$ws.Range("C2").Value = "//*[text()='Search']"

# Search text box.
$ws.Range("A3").Value = "txtSearch"
$ws.Range("B3").Value = "by_id"
$ws.Range("C3").Value = "tbSearchCriteria"

# Search button.
$ws.Range("A4").Value = "btnSearch"
$ws.Range("B4").Value = "by_xpath"
$ws.Range("C4").Value = "//*[@class='btn btn-default btn-sm']"

# Result grid row.
$ws.Range("A5").Value = "tableRow"
$ws.Range("B5").Value = "by_id"
$ws.Range("C5").Value = "gridLease"

# Column widths, sized roughly to fit the longest value in each column.
$ws.Columns.Item(1).ColumnWidth = 20.42578125
$ws.Columns.Item(2).ColumnWidth = 37.85546875
$ws.Columns.Item(3).ColumnWidth = 38.7109375
$ws.Columns.Item(4).ColumnWidth = 42.7109375

# Leave the selection on A4 and the new sheet as the active tab (mirrors the
# saved-workbook state of the authored change).
$ws.Range("A4").Select() | Out-Null
